# H17R10 BOM workbook — "Correcting some documentation issues"
#
# 1. The worksheet tab was mis-named after an older board revision
#    (H21R00); rename it to match the actual board/document (H17R10).
# 2. Two designator cells (A13 / A16) were left in a red "draft/needs
#    review" font from earlier editing. Clean them up by copying the
#    normal (black) formatting already used by the other designator
#    cells in the same column (e.g. A17) onto them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the sheet/tab name ------------------------------------------------
$ws.Name = "H17R10"

# --- Clear the leftover red formatting on A13 and A16 ----------------------
$ws.Range("A17").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Restore the normal cell selection/cursor position ----------------------
$ws.Range("B14").Select()
